$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("A2").NumberFormat = "@"
    $ws.Range("A2").Value = '2025-07-31'
    $ws.Range("A2").Style = "Normal"
    $ws.Range("C2").NumberFormat = "@"
    $ws.Range("C2").Value = 'BEMOL S/A'
    $ws.Range("C2").Style = "Normal"
    $ws.Range("D2").NumberFormat = "@"
    $ws.Range("D2").Value = '390878'
    $ws.Range("D2").Style = "Normal"
    $ws.Range("F2").NumberFormat = "@"
    $ws.Range("F2").Value = 'SUPORTE PARA NOTEBOOK HMASTON PRETO'
    $ws.Range("F2").Style = "Normal"
    $ws.Range("B2").Value = 7
    $ws.Range("E2").Value = 48587
    $ws.Range("G2").Value = -134
    $ws.Range("H2").Value = 1.26
    $ws.Range("I2").Value = 1.43

    # Row 3
    $ws.Range("A3").NumberFormat = "@"
    $ws.Range("A3").Value = '2025-08-04'
    $ws.Range("A3").Style = "Normal"
    $ws.Range("C3").NumberFormat = "@"
    $ws.Range("C3").Value = 'BEMOL S/A'
    $ws.Range("C3").Style = "Normal"
    $ws.Range("D3").NumberFormat = "@"
    $ws.Range("D3").Value = '391921'
    $ws.Range("D3").Style = "Normal"
    $ws.Range("F3").NumberFormat = "@"
    $ws.Range("F3").Value = 'FONE DE OUVIDO SEM FIO A GOLD V5.3'
    $ws.Range("F3").Style = "Normal"
    $ws.Range("B3").Value = 2
    $ws.Range("E3").Value = 10130
    $ws.Range("G3").Value = -1314
    $ws.Range("H3").Value = 1.05
    $ws.Range("I3").Value = 0.26

    # Row 4
    $ws.Range("A4").NumberFormat = "@"
    $ws.Range("A4").Value = '2025-08-07'
    $ws.Range("A4").Style = "Normal"
    $ws.Range("C4").NumberFormat = "@"
    $ws.Range("C4").Value = 'BEMOL S/A'
    $ws.Range("C4").Style = "Normal"
    $ws.Range("D4").NumberFormat = "@"
    $ws.Range("D4").Value = '393760'
    $ws.Range("D4").Style = "Normal"
    $ws.Range("F4").NumberFormat = "@"
    $ws.Range("F4").Value = 'HEADSET GAMER PLAYER PLUS LED 7 CORES 2M DRIVER 50MM PRETO LETRON'
    $ws.Range("F4").Style = "Normal"
    $ws.Range("B4").Value = 2
    $ws.Range("E4").Value = 14152
    $ws.Range("G4").Value = -13
    $ws.Range("H4").Value = 1.08
    $ws.Range("I4").Value = 0.29

    # Row 5
    $ws.Range("A5").NumberFormat = "@"
    $ws.Range("A5").Value = '2025-08-07'
    $ws.Range("A5").Style = "Normal"
    $ws.Range("C5").NumberFormat = "@"
    $ws.Range("C5").Value = 'BEMOL S/A'
    $ws.Range("C5").Style = "Normal"
    $ws.Range("D5").NumberFormat = "@"
    $ws.Range("D5").Value = '393791'
    $ws.Range("D5").Style = "Normal"
    $ws.Range("F5").NumberFormat = "@"
    $ws.Range("F5").Value = 'MINI VENTILADOR RECARREGAVEL E PORTATIL, COM PREGADOR MATERIAL ABS E COMPONENTES ELETRONICOS'
    $ws.Range("F5").Style = "Normal"
    $ws.Range("B5").Value = 3
    $ws.Range("E5").Value = 13965
    $ws.Range("G5").Value = -25
    $ws.Range("H5").Value = 1.2
    $ws.Range("I5").Value = 0.5

    # Row 6
    $ws.Range("A6").NumberFormat = "@"
    $ws.Range("A6").Value = '2025-08-08'
    $ws.Range("A6").Style = "Normal"
    $ws.Range("C6").NumberFormat = "@"
    $ws.Range("C6").Value = 'BEMOL S/A'
    $ws.Range("C6").Style = "Normal"
    $ws.Range("D6").NumberFormat = "@"
    $ws.Range("D6").Value = '394429'
    $ws.Range("D6").Style = "Normal"
    $ws.Range("F6").NumberFormat = "@"
    $ws.Range("F6").Value = 'FONE DE OUVIDO SEM FIO A GOLD V5.3'
    $ws.Range("F6").Style = "Normal"
    $ws.Range("B6").Value = 2
    $ws.Range("E6").Value = 10130
    $ws.Range("G6").Value = -1314
    $ws.Range("H6").Value = 1.05
    $ws.Range("I6").Value = 0.26

    # Row 7
    $ws.Range("A7").NumberFormat = "@"
    $ws.Range("A7").Value = '2025-08-11'
    $ws.Range("A7").Style = "Normal"
    $ws.Range("C7").NumberFormat = "@"
    $ws.Range("C7").Value = 'BEMOL S/A'
    $ws.Range("C7").Style = "Normal"
    $ws.Range("D7").NumberFormat = "@"
    $ws.Range("D7").Value = '396518'
    $ws.Range("D7").Style = "Normal"
    $ws.Range("F7").NumberFormat = "@"
    $ws.Range("F7").Value = 'FONE DE OUVIDO SEM FIO A GOLD V5.3'
    $ws.Range("F7").Style = "Normal"
    $ws.Range("B7").Value = 2
    $ws.Range("E7").Value = 10130
    $ws.Range("G7").Value = -1314
    $ws.Range("H7").Value = 1.05
    $ws.Range("I7").Value = 0.26

    # Row 8
    $ws.Range("A8").NumberFormat = "@"
    $ws.Range("A8").Value = '2025-08-11'
    $ws.Range("A8").Style = "Normal"
    $ws.Range("C8").NumberFormat = "@"
    $ws.Range("C8").Value = 'BEMOL S/A'
    $ws.Range("C8").Style = "Normal"
    $ws.Range("D8").NumberFormat = "@"
    $ws.Range("D8").Value = '396572'
    $ws.Range("D8").Style = "Normal"
    $ws.Range("F8").NumberFormat = "@"
    $ws.Range("F8").Value = 'FONE DE OUVIDO SEM FIO A GOLD V5.3'
    $ws.Range("F8").Style = "Normal"
    $ws.Range("B8").Value = 2
    $ws.Range("E8").Value = 10130
    $ws.Range("G8").Value = -1314
    $ws.Range("H8").Value = 1.05
    $ws.Range("I8").Value = 0.26

    # Row 9
    $ws.Range("A9").NumberFormat = "@"
    $ws.Range("A9").Value = '2025-08-12'
    $ws.Range("A9").Style = "Normal"
    $ws.Range("C9").NumberFormat = "@"
    $ws.Range("C9").Value = 'BEMOL S/A'
    $ws.Range("C9").Style = "Normal"
    $ws.Range("D9").NumberFormat = "@"
    $ws.Range("D9").Value = '397270'
    $ws.Range("D9").Style = "Normal"
    $ws.Range("F9").NumberFormat = "@"
    $ws.Range("F9").Value = 'MOUSE SEM FIO 3 BOTOES 1000DPI COLOR FIT AZUL 1709 R8'
    $ws.Range("F9").Style = "Normal"
    $ws.Range("B9").Value = 2
    $ws.Range("E9").Value = 13244
    $ws.Range("G9").Value = -15
    $ws.Range("H9").Value = 1.06
    $ws.Range("I9").Value = 0.25

    # Row 10
    $ws.Range("A10").NumberFormat = "@"
    $ws.Range("A10").Value = '2025-08-12'
    $ws.Range("A10").Style = "Normal"
    $ws.Range("C10").NumberFormat = "@"
    $ws.Range("C10").Value = 'BEMOL S/A'
    $ws.Range("C10").Style = "Normal"
    $ws.Range("D10").NumberFormat = "@"
    $ws.Range("D10").Value = '397270'
    $ws.Range("D10").Style = "Normal"
    $ws.Range("F10").NumberFormat = "@"
    $ws.Range("F10").Value = 'MOUSE SEM FIO 3 BOTOES 1000DPI COLOR FIT BRANCO 1709 R8'
    $ws.Range("F10").Style = "Normal"
    $ws.Range("B10").Value = 2
    $ws.Range("E10").Value = 13544
    $ws.Range("G10").Value = 3
    $ws.Range("H10").Value = 1.07
    $ws.Range("I10").Value = 0.27

    # Row 11
    $ws.Range("A11").NumberFormat = "@"
    $ws.Range("A11").Value = '2025-08-12'
    $ws.Range("A11").Style = "Normal"
    $ws.Range("C11").NumberFormat = "@"
    $ws.Range("C11").Value = 'BEMOL S/A'
    $ws.Range("C11").Style = "Normal"
    $ws.Range("D11").NumberFormat = "@"
    $ws.Range("D11").Value = '397297'
    $ws.Range("D11").Style = "Normal"
    $ws.Range("F11").NumberFormat = "@"
    $ws.Range("F11").Value = 'FONE DE OUVIDO SEM FIO BT BASIKE FON-9856'
    $ws.Range("F11").Style = "Normal"
    $ws.Range("B11").Value = 2
    $ws.Range("E11").Value = 12945
    $ws.Range("G11").Value = -91
    $ws.Range("H11").Value = 1.03
    $ws.Range("I11").Value = 0.17

    # Row 12
    $ws.Range("A12").NumberFormat = "@"
    $ws.Range("A12").Value = '2025-08-12'
    $ws.Range("A12").Style = "Normal"
    $ws.Range("C12").NumberFormat = "@"
    $ws.Range("C12").Value = 'BEMOL S/A'
    $ws.Range("C12").Style = "Normal"
    $ws.Range("D12").NumberFormat = "@"
    $ws.Range("D12").Value = '397345'
    $ws.Range("D12").Style = "Normal"
    $ws.Range("F12").NumberFormat = "@"
    $ws.Range("F12").Value = 'Kit Smartwatch Inova Com Pulseira Respiravel Preto Fone E Carregador'
    $ws.Range("F12").Style = "Normal"
    $ws.Range("B12").Value = 2
    $ws.Range("E12").Value = 396985
    $ws.Range("G12").Value = -45
    $ws.Range("H12").Value = 1.02
    $ws.Range("I12").Value = 0.15

    # Row 13
    $ws.Range("A13").NumberFormat = "@"
    $ws.Range("A13").Value = '2025-08-13'
    $ws.Range("A13").Style = "Normal"
    $ws.Range("C13").NumberFormat = "@"
    $ws.Range("C13").Value = 'BEMOL S/A'
    $ws.Range("C13").Style = "Normal"
    $ws.Range("D13").NumberFormat = "@"
    $ws.Range("D13").Value = '398131'
    $ws.Range("D13").Style = "Normal"
    $ws.Range("F13").NumberFormat = "@"
    $ws.Range("F13").Value = 'FONE DE OUVIDO SEM FIO A GOLD V5.3'
    $ws.Range("F13").Style = "Normal"
    $ws.Range("B13").Value = 2
    $ws.Range("E13").Value = 10130
    $ws.Range("G13").Value = -1314
    $ws.Range("H13").Value = 1.05
    $ws.Range("I13").Value = 0.26

    # Row 14
    $ws.Range("A14").NumberFormat = "@"
    $ws.Range("A14").Value = '2025-08-13'
    $ws.Range("A14").Style = "Normal"
    $ws.Range("C14").NumberFormat = "@"
    $ws.Range("C14").Value = 'BEMOL S/A'
    $ws.Range("C14").Style = "Normal"
    $ws.Range("D14").NumberFormat = "@"
    $ws.Range("D14").Value = '398157'
    $ws.Range("D14").Style = "Normal"
    $ws.Range("F14").NumberFormat = "@"
    $ws.Range("F14").Value = 'FONE HEADSET FONE DE OUVIDO PEI-P9 MUSIC POWER'
    $ws.Range("F14").Style = "Normal"
    $ws.Range("B14").Value = 2
    $ws.Range("E14").Value = 11436
    $ws.Range("G14").Value = 1
    $ws.Range("H14").Value = 1.02
    $ws.Range("I14").Value = 0.13

